# Update the MSME indicator figures on the "Summary" sheet with more
# precise decimal values. The original cells hold their numeric-looking
# values as plain text (General format), so we re-enter each value with a
# leading apostrophe to keep it stored as text rather than letting Excel
# auto-convert it to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises (absolute #): Micro / SMEs
$ws.Range("B10").Value = "'385826.54"
$ws.Range("C10").Value = "'14408.46"

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B11").Value = "'19.67"
$ws.Range("C11").Value = "'0.73"
$ws.Range("D11").Value = "'20.41"

# Employment (% of total): MSMEs
$ws.Range("D12").Value = "'26.39"
